# Update public EPEX spot / Gaz / CO2 price workbook with the latest day of data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column (AL) with the 21-jul prices.
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("AL1").Value = "21-jul"

# Reuse the header's style (bold, centered, bordered) for the new header cell,
# the same way the rest of row 1 / column A are styled.
$wsSpot.Range("AK1").Copy()
$wsSpot.Range("AL1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$spotValues = @(58.02, 40.53, 26.72, 21.12, 17.88, 22.77, 22.08, 25.09, 19.41, 0, 0, 6, 15.23, 22.63, 17.22, 16.01, 8.93, 16.19, 21.61, 42.98, 29.33, 34.2, 55.35, 44.9)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 38).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the two latest daily prices (rows 35-36).
#
# The date column (A) stores plain text such as "2025-07-18", not a real
# date value, like the rest of the column. Entering that text straight into
# a General formatted cell would make Excel auto-convert it to a date
# serial number, so the cell is briefly switched to Text format while the
# value is entered and then restored to the same (unstyled / General) look
# as the cells above it.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A35").NumberFormat = "@"
$wsGaz.Range("A35").Value = "2025-07-19"
$wsGaz.Range("A34").Copy()
$wsGaz.Range("A35").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$wsGaz.Range("B35").Value = 32.65

$wsGaz.Range("A36").NumberFormat = "@"
$wsGaz.Range("A36").Value = "2025-07-20"
$wsGaz.Range("A34").Copy()
$wsGaz.Range("A36").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$wsGaz.Range("B36").Value = 32.65

# ---------------------------------------------------------------------------
# Sheet "CO2": append the two latest daily prices (rows 35-36).
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A35").NumberFormat = "@"
$wsCo2.Range("A35").Value = "2025-07-19"
$wsCo2.Range("A34").Copy()
$wsCo2.Range("A35").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$wsCo2.Range("B35").Value = 69.2

$wsCo2.Range("A36").NumberFormat = "@"
$wsCo2.Range("A36").Value = "2025-07-20"
$wsCo2.Range("A34").Copy()
$wsCo2.Range("A36").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$wsCo2.Range("B36").Value = 69.2

$wb.Save()
